$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A2").Value = "Location"
$ws.Range("A3").Value = "Random"
$ws.Range("A4").Value = "Ticket Inquire"
$ws.Range("A5").Value = "Winnings and Numbers"
$ws.Range("A6").Value = "Ticket Inquire Now"
$ws.Range("A7").Value = "Next Gen Games"

$ws.Range("A8").Select()
